$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 170898984.535841
$ws.Range("D2").Value = 97.546297

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 37151728.931981
$ws.Range("D3").Value = 10.602794
$ws.Range("E3").Value = 0.000034

# Row 4 - Residuals
$ws.Range("B4").Value = 581656760.2900831
$ws.Range("C4").Value = 332

# Row 5 - SM-Control
$ws.Range("G5").Value = -608.985486
$ws.Range("H5").Value = -1088.882887
$ws.Range("I5").Value = -129.088084
$ws.Range("J5").Value = 0.00847

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = 82.161582
$ws.Range("H6").Value = -421.893392
$ws.Range("I6").Value = 586.216557
$ws.Range("J6").Value = 0.92206

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 691.147068
$ws.Range("H7").Value = 312.131393
$ws.Range("I7").Value = 1070.162744
$ws.Range("J7").Value = 0.000069
